# Update Work Week and Social Spending
# Rewrites the GDP-per-Capita "Data" values for Myanmar (years 1820-2010) with revised
# figures, and appends six new rows for years 2011-2016.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

function Set-TextValue($cell, $text) {
    # Force the cell's stored value to be a real *text* string (matching the
    # workbook convention of storing "Data" figures as shared strings) rather than
    # letting Excel auto-convert a numeric-looking value into a Number cell.
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# row number -> revised "Data" text value
$updates = @(
    @(2, "803"),
    @(52, "803"),
    @(83, "1114"),
    @(88, "902"),
    @(93, "968"),
    @(95, "1092"),
    @(98, "1312"),
    @(103, "1133"),
    @(108, "1298"),
    @(113, "1438"),
    @(118, "1336"),
    @(120, "1180"),
    @(132, "631"),
    @(133, "711"),
    @(134, "716"),
    @(135, "724"),
    @(136, "668"),
    @(137, "744"),
    @(138, "781"),
    @(139, "813"),
    @(140, "778"),
    @(141, "885"),
    @(142, "899"),
    @(143, "904"),
    @(144, "966"),
    @(145, "977"),
    @(146, "976"),
    @(147, "983"),
    @(148, "925"),
    @(149, "934"),
    @(150, "976"),
    @(151, "998"),
    @(152, "1023"),
    @(153, "1036"),
    @(154, "1023"),
    @(155, "1001"),
    @(156, "1033"),
    @(157, "1057"),
    @(158, "1103"),
    @(159, "1148"),
    @(160, "1205"),
    @(161, "1240"),
    @(162, "1320"),
    @(163, "1369"),
    @(164, "1417"),
    @(165, "1451"),
    @(166, "1492"),
    @(167, "1506"),
    @(168, "1465"),
    @(169, "1382"),
    @(170, "1207"),
    @(171, "1234"),
    @(172, "1253"),
    @(173, "1205.07146138499"),
    @(174, "1283.82548302487"),
    @(175, "1321.81019355356"),
    @(176, "1373.96288657968"),
    @(177, "1423.65048507669"),
    @(178, "1469.03173707452"),
    @(179, "1505.98046383745"),
    @(180, "1547.02833973273"),
    @(181, "1667.66122635631"),
    @(182, "1843.494517528"),
    @(183, "1996.37368678253"),
    @(184, "2175.88343276733"),
    @(185, "2410.56976205422"),
    @(186, "2664.94125386245"),
    @(187, "2947.0989157759"),
    @(188, "3246.69342175902"),
    @(189, "3542.58842689859"),
    @(190, "3585.76522166707"),
    @(191, "3673.43665669723"),
    @(192, "3772.99854122546"),
)

foreach ($pair in $updates) {
    $row = $pair[0]
    $value = $pair[1]
    Set-TextValue $ws.Cells.Item($row, 5) $value
}

# Append the new years (2011-2016)
$newRows = @(
    @(193, 2011, "3884"),
    @(194, 2012, "4124"),
    @(195, 2013, "4424"),
    @(196, 2014, "4728"),
    @(197, 2015, "5021"),
    @(198, 2016, "5284"),
)

foreach ($row3 in $newRows) {
    $r = $row3[0]
    $year = $row3[1]
    $value = $row3[2]
    $ws.Cells.Item($r, 1).Value = 104
    $ws.Cells.Item($r, 2).Value = "Myanmar"
    $ws.Cells.Item($r, 3).Value = "GDP per Capita"
    $ws.Cells.Item($r, 4).Value = $year
    Set-TextValue $ws.Cells.Item($r, 5) $value
}
